$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model name ordering (column A) for rows 2..26
$names = @{
    2  = "model_14_8_0"
    3  = "model_14_8_22"
    4  = "model_14_8_21"
    5  = "model_14_8_20"
    6  = "model_14_8_19"
    7  = "model_14_8_18"
    8  = "model_14_8_17"
    9  = "model_14_8_16"
    10 = "model_14_8_15"
    11 = "model_14_8_14"
    12 = "model_14_8_13"
    13 = "model_14_8_23"
    14 = "model_14_8_12"
    15 = "model_14_8_10"
    16 = "model_14_8_9"
    17 = "model_14_8_8"
    18 = "model_14_8_7"
    19 = "model_14_8_6"
    20 = "model_14_8_5"
    21 = "model_14_8_4"
    22 = "model_14_8_3"
    23 = "model_14_8_2"
    24 = "model_14_8_1"
    25 = "model_14_8_11"
    26 = "model_14_8_24"
}

# New metric values (columns B..Q), identical for every data row (2..26)
$values = @{
    "B" = [double]"0.999997025781996"
    "C" = [double]"0.9990399586905852"
    "D" = [double]"0.9999999999998939"
    "E" = [double]"0.9999998150907792"
    "F" = [double]"0.9999998953322146"
    "G" = [double]"2.776301953652792e-06"
    "H" = [double]"0.0008961564213956752"
    "I" = [double]"1.044290021474919e-13"
    "J" = [double]"2.548652552803534e-07"
    "K" = [double]"1.274326798546778e-07"
    "L" = [double]"0.0001004477823987822"
    "M" = [double]"0.001666223860606009"
    "N" = [double]"1.000007931248011"
    "O" = [double]"0.001737158459910895"
    "P" = [double]"91.58878149661076"
    "Q" = [double]"131.8116837172614"
}

for ($row = 2; $row -le 26; $row++) {
    $ws.Range("A$row").Value = $names[$row]
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
